$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 1288
$ws.Range("G6").Value = 68
$ws.Range("F7").Value = 63
$ws.Range("F10").Value = 441
$ws.Range("F11").Value = 811
$ws.Range("F13").Value = 742
$ws.Range("F14").Value = 306
$ws.Range("F15").Value = 458
$ws.Range("F16").Value = 90
$ws.Range("F17").Value = 1048
$ws.Range("F18").Value = 489
$ws.Range("F19").Value = 287
$ws.Range("F21").Value = 100
$ws.Range("F22").Value = 212
$ws.Range("F26").Value = 428
$ws.Range("F27").Value = 271

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 48
$ws.Range("F11").Value = 158
$ws.Range("F12").Value = 140
$ws.Range("F13").Value = 36
$ws.Range("F14").Value = 7

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 1288
$ws.Range("G8").Value = 68
$ws.Range("F10").Value = 63
$ws.Range("F13").Value = 48
$ws.Range("F17").Value = 441
$ws.Range("F18").Value = 811
$ws.Range("F20").Value = 742
$ws.Range("F21").Value = 306
$ws.Range("F22").Value = 458
$ws.Range("F23").Value = 90
$ws.Range("F24").Value = 1048
$ws.Range("F25").Value = 489
$ws.Range("F28").Value = 287
$ws.Range("F31").Value = 100
$ws.Range("F32").Value = 158
$ws.Range("F33").Value = 212
$ws.Range("F36").Value = 140
$ws.Range("F37").Value = 36
$ws.Range("F39").Value = 7
$ws.Range("F41").Value = 428
$ws.Range("F42").Value = 271

